# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.126.07'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.624.27'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.523'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0848'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '1.631.18'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.66'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').Value = '27.096.60'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '0.0₃0744'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.55%  '
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').Value = '1.345.59'
$ws.Range('E33').Value = '  +5.51%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0178'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.858'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.803'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('E41').Value = '  +6.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('D44').Value = '1.761.61'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.64'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.866'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +29.63%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.91%  '
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0993'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.05%  '
